# Auto-generated: sync of market-price / profit columns (H-N) on the
# per-job Leve tables, pulled from the scheduled price-update runner.
# Each row below updates currentAveragePrice* / LevePrice* / LeveProfit*
# cells (columns H-N) to the freshly fetched values; a couple of rows
# lose their HQ-profit cell entirely because LevePriceHQ dropped to 0.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1141
$ws.Cells.Item(98, 9).Value = 235
$ws.Cells.Item(98, 10).Value = 2500
$ws.Cells.Item(98, 11).Value = 235
$ws.Cells.Item(98, 12).Value = 2500
$ws.Cells.Item(98, 13).Value = 1263
$ws.Cells.Item(98, 14).Value = -5496
$ws.Cells.Item(122, 8).Value = 1141
$ws.Cells.Item(122, 9).Value = 235
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 705
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = 1745
$ws.Cells.Item(122, 14).Value = -12400
$ws.Cells.Item(132, 8).Value = 4622
$ws.Cells.Item(132, 9).Value = 4874.4585
$ws.Cells.Item(132, 10).Value = 3756.4285
$ws.Cells.Item(132, 11).Value = 14623.3755
$ws.Cells.Item(132, 12).Value = 11269.2855
$ws.Cells.Item(132, 13).Value = -12093.3755
$ws.Cells.Item(132, 14).Value = -16329.2855
$ws.Cells.Item(27, 8).Value = 74990
$ws.Cells.Item(27, 10).Value = 74990
$ws.Cells.Item(27, 12).Value = 74990
$ws.Cells.Item(27, 14).Value = -75374

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 430275.97
$ws.Cells.Item(32, 9).Value = 499936.3
$ws.Cells.Item(32, 11).Value = 499936.3
$ws.Cells.Item(32, 13).Value = -499649.3
$ws.Cells.Item(61, 8).Value = 3551.577
$ws.Cells.Item(61, 9).Value = 2987.7693
$ws.Cells.Item(61, 11).Value = 2987.7693
$ws.Cells.Item(61, 13).Value = -2775.7693
$ws.Cells.Item(132, 8).Value = 4303.4443
$ws.Cells.Item(132, 9).Value = 4400.278
$ws.Cells.Item(132, 10).Value = 4109.778
$ws.Cells.Item(132, 11).Value = 13200.834
$ws.Cells.Item(132, 12).Value = 12329.334
$ws.Cells.Item(132, 13).Value = -10670.834
$ws.Cells.Item(132, 14).Value = -17389.334
$ws.Cells.Item(136, 8).Value = 3551.577
$ws.Cells.Item(136, 9).Value = 2987.7693
$ws.Cells.Item(136, 11).Value = 8963.3079
$ws.Cells.Item(136, 13).Value = -6413.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(28, 8).Value = 54771
$ws.Cells.Item(28, 10).Value = 54771
$ws.Cells.Item(28, 12).Value = 54771
$ws.Cells.Item(28, 14).Value = -55359
$ws.Cells.Item(40, 8).Value = 66000
$ws.Cells.Item(40, 10).Value = 66000
$ws.Cells.Item(40, 12).Value = 66000
$ws.Cells.Item(40, 14).Value = -66530
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 146142.86
$ws.Cells.Item(35, 9).Value = 146142.86
$ws.Cells.Item(35, 11).Value = 146142.86
$ws.Cells.Item(35, 13).Value = -145848.86
$ws.Cells.Item(51, 8).Value = 9997.666999999999
$ws.Cells.Item(51, 10).Value = 9997.666999999999
$ws.Cells.Item(51, 12).Value = 9997.666999999999
$ws.Cells.Item(51, 14).Value = -11469.667
$ws.Cells.Item(59, 8).Value = 15140.714
$ws.Cells.Item(59, 10).Value = 15140.714
$ws.Cells.Item(59, 12).Value = 15140.714
$ws.Cells.Item(59, 14).Value = -17430.714
$ws.Cells.Item(60, 8).Value = 9999.799999999999
$ws.Cells.Item(60, 10).Value = 9999.799999999999
$ws.Cells.Item(60, 12).Value = 9999.799999999999
$ws.Cells.Item(60, 14).Value = -11021.8
$ws.Cells.Item(61, 8).Value = 9997.666999999999
$ws.Cells.Item(61, 10).Value = 9997.666999999999
$ws.Cells.Item(61, 12).Value = 9997.666999999999
$ws.Cells.Item(61, 14).Value = -10693.667
$ws.Cells.Item(68, 8).Value = 23749.416
$ws.Cells.Item(68, 10).Value = 23749.416
$ws.Cells.Item(68, 12).Value = 23749.416
$ws.Cells.Item(68, 14).Value = -25247.416
$ws.Cells.Item(71, 8).Value = 23749.416
$ws.Cells.Item(71, 10).Value = 23749.416
$ws.Cells.Item(71, 12).Value = 71248.24800000001
$ws.Cells.Item(71, 14).Value = -78736.24800000001
$ws.Cells.Item(74, 8).Value = 28998.223
$ws.Cells.Item(74, 10).Value = 28998.223
$ws.Cells.Item(74, 12).Value = 28998.223
$ws.Cells.Item(74, 14).Value = -30746.223
$ws.Cells.Item(77, 8).Value = 28998.223
$ws.Cells.Item(77, 10).Value = 28998.223
$ws.Cells.Item(77, 12).Value = 86994.66900000001
$ws.Cells.Item(77, 14).Value = -95730.66900000001
$ws.Cells.Item(132, 8).Value = 5954276.5
$ws.Cells.Item(132, 9).Value = 1278.4117
$ws.Cells.Item(132, 11).Value = 3835.2351
$ws.Cells.Item(132, 13).Value = -1305.2351

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 11668778
$ws.Cells.Item(4, 9).Value = 6668389
$ws.Cells.Item(4, 10).Value = 16669167
$ws.Cells.Item(4, 11).Value = 20005167
$ws.Cells.Item(4, 12).Value = 50007501
$ws.Cells.Item(4, 13).Value = -20005055
$ws.Cells.Item(4, 14).Value = -50007725
$ws.Cells.Item(39, 8).Value = 1399.091
$ws.Cells.Item(39, 10).Value = 1429.3024
$ws.Cells.Item(39, 12).Value = 4287.9072
$ws.Cells.Item(39, 14).Value = -4875.9072
$ws.Cells.Item(92, 8).Value = 849.5625
$ws.Cells.Item(92, 10).Value = 859.5333000000001
$ws.Cells.Item(92, 12).Value = 2578.5999
$ws.Cells.Item(92, 14).Value = -5074.5999
$ws.Cells.Item(110, 8).Value = 12502.143
$ws.Cells.Item(110, 10).Value = 14148.117
$ws.Cells.Item(110, 12).Value = 42444.351
$ws.Cells.Item(110, 14).Value = -50624.351
$ws.Cells.Item(122, 8).Value = 9721.182000000001
$ws.Cells.Item(122, 10).Value = 26024.5
$ws.Cells.Item(122, 12).Value = 234220.5
$ws.Cells.Item(122, 14).Value = -239120.5
$ws.Cells.Item(131, 8).Value = 861.5854
$ws.Cells.Item(131, 10).Value = 1060.8064
$ws.Cells.Item(131, 12).Value = 3182.4192
$ws.Cells.Item(131, 14).Value = -13262.4192
$ws.Cells.Item(137, 8).Value = 9144.764999999999
$ws.Cells.Item(137, 9).Value = 12533.8
$ws.Cells.Item(137, 10).Value = 4303.2856
$ws.Cells.Item(137, 11).Value = 37601.39999999999
$ws.Cells.Item(137, 12).Value = 12909.8568
$ws.Cells.Item(137, 13).Value = -32501.39999999999
$ws.Cells.Item(137, 14).Value = -23109.8568
$ws.Cells.Item(139, 8).Value = 2412.3809
$ws.Cells.Item(139, 9).Value = 2290.6155
$ws.Cells.Item(139, 10).Value = 2610.25
$ws.Cells.Item(139, 11).Value = 6871.8465
$ws.Cells.Item(139, 12).Value = 7830.75
$ws.Cells.Item(139, 13).Value = -1731.8465
$ws.Cells.Item(139, 14).Value = -18110.75
$ws.Cells.Item(140, 8).Value = 1203.8889
$ws.Cells.Item(140, 9).Value = 923.6
$ws.Cells.Item(140, 10).Value = 2004.7142
$ws.Cells.Item(140, 11).Value = 2770.8
$ws.Cells.Item(140, 12).Value = 6014.142599999999
$ws.Cells.Item(140, 13).Value = 2409.2
$ws.Cells.Item(140, 14).Value = -16374.1426
$ws.Cells.Item(141, 8).Value = 5369.222
$ws.Cells.Item(141, 9).Value = 4244.1177
$ws.Cells.Item(141, 11).Value = 12732.3531
$ws.Cells.Item(141, 13).Value = -7552.3531

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5687.2666
$ws.Cells.Item(70, 9).Value = 5591.706
$ws.Cells.Item(70, 10).Value = 5982.636
$ws.Cells.Item(70, 11).Value = 5591.706
$ws.Cells.Item(70, 12).Value = 5982.636
$ws.Cells.Item(70, 13).Value = -5321.706
$ws.Cells.Item(70, 14).Value = -6522.636
$ws.Cells.Item(73, 8).Value = 5687.2666
$ws.Cells.Item(73, 9).Value = 5591.706
$ws.Cells.Item(73, 10).Value = 5982.636
$ws.Cells.Item(73, 11).Value = 5591.706
$ws.Cells.Item(73, 12).Value = 5982.636
$ws.Cells.Item(73, 13).Value = -4655.706
$ws.Cells.Item(73, 14).Value = -7854.636
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()

